$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 17.41485066666667
$ws.Range("H2").Value = 52.24455200000001
$ws.Range("I2").Value = 0.1047285618770465
$ws.Range("J2").Value = 0.1047285618770465
$ws.Range("M2").Value = 2.33201
$ws.Range("N2").Value = 6.99603
$ws.Range("O2").Value = 0.3303263034789547
$ws.Range("P2").Value = 0.3303263034789548
$ws.Range("Q2").Value = 40.61160590317334
$ws.Range("R2").Value = 365.50445312856
$ws.Range("S2").Value = 0.03459459871351176
$ws.Range("T2").Value = 0.03459459871351176

$ws.Range("G3").Value = 17.41485066666667
$ws.Range("H3").Value = 52.24455200000001
$ws.Range("I3").Value = 0.1047285618770465
$ws.Range("J3").Value = 0.1047285618770465
$ws.Range("N3").Value = 5.238131999999999
$ws.Range("O3").Value = 0.2473249515360603
$ws.Range("P3").Value = 0.2473249515360603
$ws.Range("Q3").Value = 30.40709551742933
$ws.Range("R3").Value = 273.663859656864
$ws.Range("S3").Value = 0.02590198649068182
$ws.Range("T3").Value = 0.02590198649068182

$ws.Range("G4").Value = 17.41485066666667
$ws.Range("H4").Value = 52.24455200000001
$ws.Range("I4").Value = 0.1047285618770465
$ws.Range("J4").Value = 0.1047285618770465
$ws.Range("M4").Value = 1.145780666666667
$ws.Range("N4").Value = 3.437342
$ws.Range("O4").Value = 0.1622984001859565
$ws.Range("P4").Value = 0.1622984001859565
$ws.Range("Q4").Value = 19.95359920675378
$ws.Range("R4").Value = 179.582392860784
$ws.Range("S4").Value = 0.0169972780464206
$ws.Range("T4").Value = 0.0169972780464206

$ws.Range("G5").Value = 17.41485066666667
$ws.Range("H5").Value = 52.24455200000001
$ws.Range("I5").Value = 0.1047285618770465
$ws.Range("J5").Value = 0.1047285618770465
$ws.Range("M5").Value = 1.835881666666667
$ws.Range("N5").Value = 5.507645
$ws.Range("O5").Value = 0.2600503447990285
$ws.Range("P5").Value = 0.2600503447990286
$ws.Range("Q5").Value = 31.97160506667112
$ws.Range("R5").Value = 287.74444560004
$ws.Range("S5").Value = 0.02723469862643234
$ws.Range("T5").Value = 0.02723469862643235

$ws.Range("I6").Value = 0.1785014126970782
$ws.Range("J6").Value = 0.1785014126970782
$ws.Range("M6").Value = 2.33201
$ws.Range("N6").Value = 6.99603
$ws.Range("O6").Value = 0.3303263034789547
$ws.Range("P6").Value = 0.3303263034789548
$ws.Range("Q6").Value = 69.21921676079333
$ws.Range("R6").Value = 622.9729508471401
$ws.Range("S6").Value = 0.05896371182199718
$ws.Range("T6").Value = 0.0589637118219972

$ws.Range("I7").Value = 0.1785014126970782
$ws.Range("J7").Value = 0.1785014126970782
$ws.Range("N7").Value = 5.238131999999999
$ws.Range("O7").Value = 0.2473249515360603
$ws.Range("P7").Value = 0.2473249515360603
$ws.Range("Q7").Value = 51.82644933335732
$ws.Range("R7").Value = 466.4380440002159
$ws.Range("S7").Value = 0.04414785324442315
$ws.Range("T7").Value = 0.04414785324442316

$ws.Range("I8").Value = 0.1785014126970782
$ws.Range("J8").Value = 0.1785014126970782
$ws.Range("M8").Value = 1.145780666666667
$ws.Range("N8").Value = 3.437342
$ws.Range("O8").Value = 0.1622984001859565
$ws.Range("P8").Value = 0.1622984001859565
$ws.Range("Q8").Value = 34.00930541735511
$ws.Range("R8").Value = 306.083748756196
$ws.Range("S8").Value = 0.02897049371166896
$ws.Range("T8").Value = 0.02897049371166897

$ws.Range("I9").Value = 0.1785014126970782
$ws.Range("J9").Value = 0.1785014126970782
$ws.Range("M9").Value = 1.835881666666667
$ws.Range("N9").Value = 5.507645
$ws.Range("O9").Value = 0.2600503447990285
$ws.Range("P9").Value = 0.2600503447990286
$ws.Range("Q9").Value = 54.49303006083445
$ws.Range("R9").Value = 490.43727054751
$ws.Range("S9").Value = 0.04641935391898886
$ws.Range("T9").Value = 0.04641935391898888

$ws.Range("G10").Value = 84.03051233333333
$ws.Range("H10").Value = 252.091537
$ws.Range("I10").Value = 0.5053385113032314
$ws.Range("J10").Value = 0.5053385113032314
$ws.Range("M10").Value = 2.33201
$ws.Range("N10").Value = 6.99603
$ws.Range("O10").Value = 0.3303263034789547
$ws.Range("P10").Value = 0.3303263034789548
$ws.Range("Q10").Value = 195.9599950664567
$ws.Range("R10").Value = 1763.63995559811
$ws.Range("S10").Value = 0.1669266024443544
$ws.Range("T10").Value = 0.1669266024443544

$ws.Range("G11").Value = 84.03051233333333
$ws.Range("H11").Value = 252.091537
$ws.Range("I11").Value = 0.5053385113032314
$ws.Range("J11").Value = 0.5053385113032314
$ws.Range("N11").Value = 5.238131999999999
$ws.Range("O11").Value = 0.2473249515360603
$ws.Range("P11").Value = 0.2473249515360603
$ws.Range("Q11").Value = 146.7209718765426
$ws.Range("R11").Value = 1320.488746888884
$ws.Range("S11").Value = 0.1249828228173765
$ws.Range("T11").Value = 0.1249828228173766

$ws.Range("G12").Value = 84.03051233333333
$ws.Range("H12").Value = 252.091537
$ws.Range("I12").Value = 0.5053385113032314
$ws.Range("J12").Value = 0.5053385113032314
$ws.Range("M12").Value = 1.145780666666667
$ws.Range("N12").Value = 3.437342
$ws.Range("O12").Value = 0.1622984001859565
$ws.Range("P12").Value = 0.1622984001859565
$ws.Range("Q12").Value = 96.28053644162821
$ws.Range("R12").Value = 866.5248279746539
$ws.Range("S12").Value = 0.08201563193686733
$ws.Range("T12").Value = 0.08201563193686734

$ws.Range("G13").Value = 84.03051233333333
$ws.Range("H13").Value = 252.091537
$ws.Range("I13").Value = 0.5053385113032314
$ws.Range("J13").Value = 0.5053385113032314
$ws.Range("M13").Value = 1.835881666666667
$ws.Range("N13").Value = 5.507645
$ws.Range("O13").Value = 0.2600503447990285
$ws.Range("P13").Value = 0.2600503447990286
$ws.Range("Q13").Value = 154.2700770333739
$ws.Range("R13").Value = 1388.430693300365
$ws.Range("S13").Value = 0.1314134541046331
$ws.Range("T13").Value = 0.1314134541046331

$ws.Range("G14").Value = 35.158014
$ws.Range("H14").Value = 105.474042
$ws.Range("I14").Value = 0.2114315141226439
$ws.Range("J14").Value = 0.2114315141226439
$ws.Range("M14").Value = 2.33201
$ws.Range("N14").Value = 6.99603
$ws.Range("O14").Value = 0.3303263034789547
$ws.Range("P14").Value = 0.3303263034789548
$ws.Range("Q14").Value = 81.98884022814001
$ws.Range("R14").Value = 737.89956205326
$ws.Range("S14").Value = 0.06984139049909135
$ws.Range("T14").Value = 0.06984139049909137

$ws.Range("G15").Value = 35.158014
$ws.Range("H15").Value = 105.474042
$ws.Range("I15").Value = 0.2114315141226439
$ws.Range("J15").Value = 0.2114315141226439
$ws.Range("N15").Value = 5.238131999999999
$ws.Range("O15").Value = 0.2473249515360603
$ws.Range("P15").Value = 0.2473249515360603
$ws.Range("Q15").Value = 61.387439396616
$ws.Range("R15").Value = 552.486954569544
$ws.Range("S15").Value = 0.05229228898357874
$ws.Range("T15").Value = 0.05229228898357875

$ws.Range("G16").Value = 35.158014
$ws.Range("H16").Value = 105.474042
$ws.Range("I16").Value = 0.2114315141226439
$ws.Range("J16").Value = 0.2114315141226439
$ws.Range("M16").Value = 1.145780666666667
$ws.Range("N16").Value = 3.437342
$ws.Range("O16").Value = 0.1622984001859565
$ws.Range("P16").Value = 0.1622984001859565
$ws.Range("Q16").Value = 40.283372719596
$ws.Range("R16").Value = 362.550354476364
$ws.Range("S16").Value = 0.03431499649099956
$ws.Range("T16").Value = 0.03431499649099957

$ws.Range("G17").Value = 35.158014
$ws.Range("H17").Value = 105.474042
$ws.Range("I17").Value = 0.2114315141226439
$ws.Range("J17").Value = 0.2114315141226439
$ws.Range("M17").Value = 1.835881666666667
$ws.Range("N17").Value = 5.507645
$ws.Range("O17").Value = 0.2600503447990285
$ws.Range("P17").Value = 0.2600503447990286
$ws.Range("Q17").Value = 64.54595333901
$ws.Range("R17").Value = 580.91358005109
$ws.Range("S17").Value = 0.05498283814897421
$ws.Range("T17").Value = 0.05498283814897421
